$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Field data" (sheet1.xml)
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Field data")

# Insert a new column before J ("Coordinate Method"), shifting J:W -> K:X
$ws1.Range("J1").EntireColumn.Insert()
$ws1.Range("J1").Value = "Coordinate Method"

# After the shift, the old "depth (cm)"/"depth (m)" headers land on M1/N1;
# relabel them per the new template wording.
$ws1.Range("M1").Value = "Depth (cm)"
$ws1.Range("N1").Value = "Depth m)"

# Transect Type column: rename the two transect categories.
$ws1.Range("F2:F37").Value = "Parallel"
$ws1.Range("F38:F70").Value = "Perpendicular"

# Quadrat IDs for the "Random" transect: zero-pad R1..R9 -> R01..R09.
$ws1.Range("G71").Value = "R01"
$ws1.Range("G72").Value = "R02"
$ws1.Range("G73").Value = "R03"
$ws1.Range("G74").Value = "R04"
$ws1.Range("G75").Value = "R05"
$ws1.Range("G76").Value = "R06"
$ws1.Range("G77").Value = "R07"
$ws1.Range("G78").Value = "R08"
$ws1.Range("G79").Value = "R09"

# ----------------------------------------------------------------------
# Sheet "ReadMe" (sheet3.xml)
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("ReadMe")

# Stray leftover "Species %" cells in row 2 are removed.
$ws3.Range("U2").ClearContents()
$ws3.Range("V2").ClearContents()

# Clarify the Quadrat definition wording.
$ws3.Range("B8").Value = "The quadrat ID, constructed using the combination of transect type, row, column, and number of metres, according to the protocol description"

# Insert a new row before row 11 describing the new "Coordinate method" field.
$ws3.Range("A11").EntireRow.Insert()
$ws3.Range("A11").Value = "Coordinate method"
$ws3.Range("B11").Value = 'Indicate whether the coordinate was obtained in the field with the GPS, or after field collection from the photo. One of "GPS" or "photo"'
